$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.375
$ws.Range("P2").Value = 0.125
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.25
$ws.Range("J3").Value = 0.1111111111111111
$ws.Range("P3").Value = 0.7777777777777778
$ws.Range("S3").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("R6").Value = 0.04761904761904762
$ws.Range("S6").Value = 0.6666666666666666
$ws.Range("B7").Value = 0.1904761904761905
$ws.Range("D7").Value = 0.04761904761904762
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("O7").Value = 0.04761904761904762
$ws.Range("Q7").Value = 0.09523809523809523
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.5238095238095238
$ws.Range("B8").Value = 0.1
$ws.Range("E8").Value = 0.01666666666666667
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.1166666666666667
$ws.Range("Q8").Value = 0.2166666666666667
$ws.Range("R8").Value = 0.05
$ws.Range("S8").Value = 0.4333333333333333
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("Q9").Value = 0.1904761904761905
$ws.Range("S9").Value = 0.6190476190476191
$ws.Range("B10").Value = 0.02962962962962963
$ws.Range("D10").Value = 0.03703703703703703
$ws.Range("F10").Value = 0.06666666666666667
$ws.Range("J10").Value = 0.1185185185185185
$ws.Range("O10").Value = 0.007407407407407408
$ws.Range("Q10").Value = 0.1925925925925926
$ws.Range("R10").Value = 0.05925925925925926
$ws.Range("S10").Value = 0.4888888888888889
$ws.Range("G11").Value = 0.1388888888888889
$ws.Range("J11").Value = 0.05555555555555555
$ws.Range("S11").Value = 0.02777777777777778
$ws.Range("G12").Value = 0.5714285714285714
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("S12").Value = 0.09523809523809523
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.1428571428571428
$ws.Range("H15").Value = 0.2105263157894737
$ws.Range("I15").Value = 0.05263157894736842
$ws.Range("J15").Value = 0.4210526315789473
$ws.Range("M15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.2631578947368421
$ws.Range("F16").Value = 0.06666666666666667
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.06666666666666667
$ws.Range("J16").Value = 0.2666666666666667
$ws.Range("K16").Value = 0.1333333333333333
$ws.Range("O16").Value = 0.1333333333333333
$ws.Range("S16").Value = 0.1333333333333333
$ws.Range("H17").Value = 0.2272727272727273
$ws.Range("I17").Value = 0.04545454545454546
$ws.Range("J17").Value = 0.4090909090909091
$ws.Range("K17").Value = 0.1363636363636364
$ws.Range("S17").Value = 0.1818181818181818
$ws.Range("H18").Value = 0.3846153846153846
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.3846153846153846
$ws.Range("O18").Value = 0.1538461538461539
$ws.Range("F19").Value = 0.01666666666666667
$ws.Range("H19").Value = 0.2055555555555555
$ws.Range("I19").Value = 0.1
$ws.Range("J19").Value = 0.35
$ws.Range("K19").Value = 0.1
$ws.Range("M19").Value = 0.03888888888888889
$ws.Range("O19").Value = 0.05555555555555555
$ws.Range("S19").Value = 0.1333333333333333

$wb.Save()
